$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'29.624.70"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +2.43%  "
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'1.860.30"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.56%  "
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'0.9992"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.06%  "
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'244.98"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +1.89%  "
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'0.6963"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +1.09%  "
$ws.Range('E6').Style = 'Normal'

$ws.Range('D7').Value = "'0.9999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +0.07%  "
$ws.Range('E7').Style = 'Normal'

$ws.Range('D8').Value = "'0.07707"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.70%  "
$ws.Range('E8').Style = 'Normal'

$ws.Range('E9').Value = "'  +0.44%  "
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'23.74"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.85%  "
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'0.07761"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  -0.53%  "
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').Value = "'  +1.77%  "
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'1.855.62"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +1.46%  "
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'92.08"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.79%  "
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'0.6930"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +2.50%  "
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = "'6.573"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +2.17%  "
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'29.611.73"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +2.41%  "
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'0.000008309"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.52%  "
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'2.101.19"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +1.24%  "
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'241.17"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  -0.67%  "
$ws.Range('E20').Style = 'Normal'

$ws.Range('E21').Value = "'  +0.87%  "
$ws.Range('E21').Style = 'Normal'

$ws.Range('D22').Value = "'0.9998"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.05%  "
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'7.604"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  +2.48%  "
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'0.9999"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.12%  "
$ws.Range('E24').Style = 'Normal'

$ws.Range('D25').Value = "'0.1501"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +1.86%  "
$ws.Range('E25').Style = 'Normal'

$ws.Range('D26').Value = "'8.926"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  +1.75%  "
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = "'159.85"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -1.00%  "
$ws.Range('E27').Style = 'Normal'

$ws.Range('D28').Value = "'18.30"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.65%  "
$ws.Range('E28').Style = 'Normal'

$ws.Range('E29').Value = "'  -0.10%  "
$ws.Range('E29').Style = 'Normal'

$ws.Range('E30').Value = "'  +1.01%  "
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'4.195"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +1.79%  "
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'1.200"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.70%  "
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'0.05089"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -0.36%  "
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'0.7739"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +2.87%  "
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Value = "'1.897"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +3.70%  "
$ws.Range('E35').Style = 'Normal'

$ws.Range('D36').Value = "'1.154"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +0.83%  "
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = "'2.683"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  +0.40%  "
$ws.Range('E37').Style = 'Normal'

$ws.Range('D38').Value = "'1.333.53"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +9.25%  "
$ws.Range('E38').Style = 'Normal'

$ws.Range('D39').Value = "'0.01873"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.60%  "
$ws.Range('E39').Style = 'Normal'

$ws.Range('D40').Value = "'2.731"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.87%  "
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Value = "'0.9782"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +6.76%  "
$ws.Range('E41').Style = 'Normal'

$ws.Range('D42').Value = "'106.63"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  -1.53%  "
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = "'5.822"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +6.46%  "
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'0.9994"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +0.09%  "
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'0.00000000127"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +4.75%  "
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'9.780"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  +2.61%  "
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'1.999.52"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.25%  "
$ws.Range('E47').Style = 'Normal'

$ws.Range('D49').Value = "'1.780"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +2.54%  "
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = "'63.71"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.42%  "
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = "'6.969"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +1.10%  "
$ws.Range('E51').Style = 'Normal'
